# Atualizado por script em 02-12-2023 20:45
#
# This script applies the diff between the pre-edit and post-edit versions
# of the laliga2 2023-2024 match results sheet:
#   1) A number of existing match rows had their F:V (match data) content
#      swapped between two adjacent rows sharing the same kickoff date
#      (column E) - i.e. the two matches that happened on the same date
#      were simply re-ordered relative to each other.
#   2) Three brand-new match rows were appended at the end of the sheet
#      (rows 190-192), extending the data range from A1:V189 to A1:V192.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap the F:V contents of the following row pairs
# ---------------------------------------------------------------------
$swapPairs = @(
    @(27,28),
    @(72,73),
    @(86,87),
    @(97,98),
    @(118,119),
    @(120,121),
    @(124,125),
    @(130,131),
    @(141,142),
    @(158,159),
    @(169,170),
    @(173,174),
    @(175,176),
    @(184,185)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("F$r1`:V$r1")
    $rng2 = $ws.Range("F$r2`:V$r2")

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value = $vals2
    $rng2.Value = $vals1
}

# ---------------------------------------------------------------------
# 2) Append three new rows (190, 191, 192) with formatting copied from
#    the last existing row (189)
# ---------------------------------------------------------------------
$ws.Range("A189:V189").Copy()
$ws.Range("A190:V192").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @{
        Row=190; Indice=189; Data=45262.67708333334
        Home="FC Cartagena SAD"; HomeGols=1; Away="Gijon"; AwayGols=0
        HomeOpenOdds=2.93; HomeOpenData="27/11/2023 21:12"
        HomeCloseOdds=3.2;  HomeCloseData="02/12/2023 16:07"
        DrawOpenOdds=3.23;  DrawOpenData="27/11/2023 21:12"
        DrawCloseOdds=3.16; DrawCloseData="02/12/2023 16:08"
        AwayOpenOdds=2.58;  AwayOpenData="27/11/2023 21:12"
        AwayCloseOdds=2.51; AwayCloseData="02/12/2023 16:08"
        Url="https://www.betexplorer.com/football/spain/laliga2/fc-cartagena-sad-gijon/88Cn0Qcd/"
    },
    @{
        Row=191; Indice=190; Data=45262.67708333334
        Home="Eldense"; HomeGols=0; Away="Tenerife"; AwayGols=3
        HomeOpenOdds=2.62; HomeOpenData="27/11/2023 21:12"
        HomeCloseOdds=2.43; HomeCloseData="02/12/2023 16:02"
        DrawOpenOdds=3.01; DrawOpenData="27/11/2023 21:12"
        DrawCloseOdds=2.93; DrawCloseData="02/12/2023 15:52"
        AwayOpenOdds=3.07; AwayOpenData="27/11/2023 21:12"
        AwayCloseOdds=3.63; AwayCloseData="02/12/2023 16:02"
        Url="https://www.betexplorer.com/football/spain/laliga2/eldense-tenerife/MB0UO9zc/"
    },
    @{
        Row=192; Indice=191; Data=45262.77083333334
        Home="Zaragoza"; HomeGols=1; Away="Leganes"; AwayGols=0
        HomeOpenOdds=2.3; HomeOpenData="26/11/2023 14:13"
        HomeCloseOdds=2.57; HomeCloseData="02/12/2023 18:29"
        DrawOpenOdds=3.08; DrawOpenData="26/11/2023 14:13"
        DrawCloseOdds=2.85; DrawCloseData="02/12/2023 18:29"
        AwayOpenOdds=3.56; AwayOpenData="26/11/2023 14:13"
        AwayCloseOdds=3.48; AwayCloseData="02/12/2023 18:29"
        Url="https://www.betexplorer.com/football/spain/laliga2/zaragoza-leganes/lSTHF5cq/"
    }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r,1).Value = $nr.Indice
    $ws.Cells.Item($r,2).Value = "spain"
    $ws.Cells.Item($r,3).Value = "laliga2"
    $ws.Cells.Item($r,4).Value = "2023-2024"
    $ws.Cells.Item($r,5).Value = $nr.Data
    $ws.Cells.Item($r,6).Value = $nr.Home
    $ws.Cells.Item($r,7).Value = $nr.HomeGols
    $ws.Cells.Item($r,8).Value = $nr.Away
    $ws.Cells.Item($r,9).Value = $nr.AwayGols
    $ws.Cells.Item($r,10).Value = $nr.HomeOpenOdds
    $ws.Cells.Item($r,11).Value = $nr.HomeOpenData
    $ws.Cells.Item($r,12).Value = $nr.HomeCloseOdds
    $ws.Cells.Item($r,13).Value = $nr.HomeCloseData
    $ws.Cells.Item($r,14).Value = $nr.DrawOpenOdds
    $ws.Cells.Item($r,15).Value = $nr.DrawOpenData
    $ws.Cells.Item($r,16).Value = $nr.DrawCloseOdds
    $ws.Cells.Item($r,17).Value = $nr.DrawCloseData
    $ws.Cells.Item($r,18).Value = $nr.AwayOpenOdds
    $ws.Cells.Item($r,19).Value = $nr.AwayOpenData
    $ws.Cells.Item($r,20).Value = $nr.AwayCloseOdds
    $ws.Cells.Item($r,21).Value = $nr.AwayCloseData
    $ws.Cells.Item($r,22).Value = $nr.Url
}

Write-Host "Edit complete"
